# Add 2022-Q3 data:
#  - "总计" sheet: insert a new row 2 with the 2022-Q3 summary figures
#  - new "2022-Q3" worksheet (cloned from "2022-Q2") inserted right after "总计"
#    with the per-fund holdings for the quarter

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet — insert the new 2022-Q3 row at the top of the data block
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)
$totals.Rows("2:2").Insert()

# the freshly inserted row picks up formatting from the row below it; reset
# it to the plain (unstyled) look used by every other data row first
$totals.Range("A3:D3").Copy()
$totals.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 3.22

# the index column is a simple 0-based row counter; every row below the
# inserted one needs to be bumped by one
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5
$totals.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet — clone "2022-Q2" (position 2) so styles/columns
#    match the sibling quarter sheets exactly, then overwrite the content
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# source sheet has 5 fund rows (2..6), the new quarter only has 3 (2..4)
$q3.Rows("5:6").Delete()

# columns B (fund code) and D:G (scale/position figures) must stay text,
# exactly like the other quarter sheets, even though they look numeric
$textRange = $q3.Range("B2:G4")
$textRange.NumberFormat = "@"

$q3.Range("B2").Value = "011056"
$q3.Range("C2").Value = "博时汇兴回报一年持有期灵活配置混合"
$q3.Range("D2").Value = "96.44"
$q3.Range("E2").Value = "52.92"
$q3.Range("F2").Value = "3.26"
$q3.Range("G2").Value = "3.1439"
$q3.Range("H2").Value = 6

$q3.Range("B3").Value = "014158"
$q3.Range("C3").Value = "博时浦惠一年持有期混合A"
$q3.Range("D3").Value = "3.87"
$q3.Range("E3").Value = "48.81"
$q3.Range("F3").Value = "1.86"
$q3.Range("G3").Value = "0.0720"
$q3.Range("H3").Value = 10

$q3.Range("B4").Value = "014159"
$q3.Range("C4").Value = "博时浦惠一年持有期混合C"
$q3.Range("D4").Value = "0.36"
$q3.Range("E4").Value = "48.81"
$q3.Range("F4").Value = "1.86"
$q3.Range("G4").Value = "0.0067"
$q3.Range("H4").Value = 10

# forcing text via NumberFormat leaves a "@" style on the cells; restore the
# plain/default look (A1 on this sheet is untouched, default-styled) so the
# new rows match the rest of the workbook
$q3.Range("A1").Copy()
$textRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totals.Activate()
